$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44329
$ws.Range("M2").Value = 56
$ws.Range("N2").Value = 9000
$ws.Range("O2").Value = 9000
$ws.Range("P2").Value = 9000
$ws.Range("R2").Value = "Región Metropolitana"
$ws.Range("S2").Value = 900

# Row 3
$ws.Range("D3").Value = 44329
$ws.Range("L3").Value = "Segunda"
$ws.Range("M3").Value = 50
$ws.Range("N3").Value = 8000
$ws.Range("O3").Value = 8000
$ws.Range("P3").Value = 8000
$ws.Range("R3").Value = "Región Metropolitana"
$ws.Range("S3").Value = 800

# Row 4
$ws.Range("D4").Value = 44302
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 45
$ws.Range("N4").Value = 10000
$ws.Range("O4").Value = 10000
$ws.Range("P4").Value = 10000
$ws.Range("S4").Value = 1000

# Row 5
$ws.Range("D5").Value = 44309
$ws.Range("M5").Value = 45
$ws.Range("N5").Value = 10000
$ws.Range("O5").Value = 10000
$ws.Range("P5").Value = 10000
$ws.Range("S5").Value = 1000

# Row 6
$ws.Range("D6").Value = 44323
$ws.Range("M6").Value = 60
$ws.Range("N6").Value = 10000
$ws.Range("O6").Value = 10000
$ws.Range("P6").Value = 10000
$ws.Range("R6").Value = "Provincia de Quillota"
$ws.Range("S6").Value = 1000

# Row 7
$ws.Range("D7").Value = 44323
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 9000
$ws.Range("P7").Value = 9000
$ws.Range("R7").Value = "Provincia de Quillota"
$ws.Range("S7").Value = 900

# Row 8
$ws.Range("D8").Value = 44308

# Row 9
$ws.Range("D9").Value = 44308
$ws.Range("L9").Value = "Segunda"
$ws.Range("M9").Value = 48
$ws.Range("N9").Value = 8000
$ws.Range("O9").Value = 8000
$ws.Range("P9").Value = 8000
$ws.Range("S9").Value = 800

# Row 10
$ws.Range("D10").Value = 44322
$ws.Range("L10").Value = "Primera"
$ws.Range("M10").Value = 56
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 10000
$ws.Range("P10").Value = 10000
$ws.Range("S10").Value = 1000

# Row 11
$ws.Range("L11").Value = "Segunda"
$ws.Range("M11").Value = 40
$ws.Range("N11").Value = 8000
$ws.Range("O11").Value = 8000
$ws.Range("P11").Value = 8000
$ws.Range("S11").Value = 800

# Row 12
$ws.Range("D12").Value = 44315
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 45
$ws.Range("N12").Value = 10000
$ws.Range("O12").Value = 10000
$ws.Range("P12").Value = 10000
$ws.Range("S12").Value = 1000

# Row 13
$ws.Range("D13").Value = 44343
$ws.Range("M13").Value = 47
$ws.Range("R13").Value = "Región Metropolitana"

# Row 14
$ws.Range("D14").Value = 44343
$ws.Range("M14").Value = 50
$ws.Range("R14").Value = "Región Metropolitana"

# Row 15
$ws.Range("D15").Value = 44343
$ws.Range("M15").Value = 58
$ws.Range("R15").Value = "Región Metropolitana"

# Row 16
$ws.Range("D16").Value = 44319
$ws.Range("M16").Value = 68
$ws.Range("N16").Value = 10000
$ws.Range("O16").Value = 10000
$ws.Range("P16").Value = 10000
$ws.Range("S16").Value = 1000

# Row 17
$ws.Range("D17").Value = 44319
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 57
$ws.Range("N17").Value = 8000
$ws.Range("O17").Value = 8000
$ws.Range("P17").Value = 8000
$ws.Range("S17").Value = 800

# Row 18
$ws.Range("D18").Value = 44307
$ws.Range("M18").Value = 40
$ws.Range("N18").Value = 10000
$ws.Range("O18").Value = 10000
$ws.Range("P18").Value = 10000
$ws.Range("S18").Value = 1000

# Row 19
$ws.Range("D19").Value = 44306
$ws.Range("L19").Value = "Primera"
$ws.Range("M19").Value = 45
$ws.Range("N19").Value = 10000
$ws.Range("O19").Value = 10000
$ws.Range("P19").Value = 10000
$ws.Range("S19").Value = 1000

# Row 20
$ws.Range("D20").Value = 44301

# Row 21
$ws.Range("D21").Value = 44328
$ws.Range("L21").Value = "Primera"
$ws.Range("M21").Value = 45
$ws.Range("N21").Value = 8000
$ws.Range("O21").Value = 8000
$ws.Range("P21").Value = 8000
$ws.Range("R21").Value = "Provincia de Quillota"
$ws.Range("S21").Value = 800

# Row 22
$ws.Range("D22").Value = 44328
$ws.Range("L22").Value = "Segunda"
$ws.Range("M22").Value = 48
$ws.Range("N22").Value = 7000
$ws.Range("O22").Value = 7000
$ws.Range("P22").Value = 7000
$ws.Range("R22").Value = "Provincia de Quillota"
$ws.Range("S22").Value = 700

# Row 23
$ws.Range("D23").Value = 44321
$ws.Range("L23").Value = "Primera"
$ws.Range("N23").Value = 9000
$ws.Range("O23").Value = 9000
$ws.Range("P23").Value = 9000
$ws.Range("R23").Value = "Provincia de Quillota"
$ws.Range("S23").Value = 900

# Row 24
$ws.Range("D24").Value = 44333
$ws.Range("L24").Value = "Especial"
$ws.Range("M24").Value = 58

# Row 25
$ws.Range("D25").Value = 44333
$ws.Range("M25").Value = 65
$ws.Range("N25").Value = 9000
$ws.Range("O25").Value = 9000
$ws.Range("P25").Value = 9000
$ws.Range("S25").Value = 900

# Row 26
$ws.Range("D26").Value = 44333
$ws.Range("M26").Value = 60

# Row 27
$ws.Range("D27").Value = 44314
$ws.Range("M27").Value = 47
$ws.Range("N27").Value = 9000
$ws.Range("O27").Value = 9000
$ws.Range("P27").Value = 9000
$ws.Range("S27").Value = 900

# Row 28
$ws.Range("D28").Value = 44312
$ws.Range("M28").Value = 48

# Row 29
$ws.Range("D29").Value = 44326
$ws.Range("M29").Value = 65

# Row 30
$ws.Range("D30").Value = 44326
$ws.Range("M30").Value = 67
$ws.Range("N30").Value = 8000
$ws.Range("O30").Value = 8000
$ws.Range("P30").Value = 8000
$ws.Range("S30").Value = 800
